$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Widen column D (value nudged so stored OOXML width rounds to exactly 45)
$ws.Columns.Item(4).ColumnWidth = 44.16666666666667

# Copy formatting from the last existing row (18) down onto the two new rows
# so the new cells get the same styles (s="8"/"12"/"7") as the rest of column A/B/C.
$ws.Range("A18:C18").Copy() | Out-Null
$ws.Range("A19:C19").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$ws.Range("A18:C18").Copy() | Out-Null
$ws.Range("A20:C20").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$ws.Application.CutCopyMode = $false

# Add two new task rows
$ws.Range("A19").Value2 = "Começar o artigo científico"
$ws.Range("B19").Value2 = "Douglas"
$ws.Range("C19").Value2 = 43257

$ws.Range("A20").Value2 = "Criar tela de tutorial"
$ws.Range("B20").Value2 = "Douglas"
$ws.Range("C20").Value2 = 43250

# Update view: scrolled position and current selection (mirrors the user's
# final on-screen state after editing B20:C20)
$ws.Application.ActiveWindow.ScrollRow = 11
$ws.Range("B20:C20").Select()
